$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values for A2:A251 (column "z_sim" data for SimCase5_zsim_SimRun1)
$values = @(1,3,3,3,3,3,3,1,2,3,3,3,1,1,1,1,3,3,3,3,3,3,2,2,3,1,3,2,1,3,3,1,3,3,1,1,3,1,3,2,3,3,2,1,3,1,3,1,3,3,3,3,1,3,1,1,3,1,1,3,2,3,1,3,3,1,3,2,2,1,3,1,1,3,1,3,1,3,3,3,3,3,3,3,2,1,1,2,3,3,3,1,1,2,3,1,2,3,1,3,1,3,2,1,2,3,3,1,1,1,3,1,1,3,2,3,3,1,3,3,2,2,1,2,1,3,1,1,3,3,3,3,2,3,2,1,3,3,3,1,3,3,3,3,3,3,2,3,1,3,3,1,1,1,3,3,2,1,3,3,3,2,3,1,3,2,3,2,3,3,2,1,1,2,3,1,1,3,3,2,1,3,1,3,1,1,1,3,3,3,3,3,3,1,1,3,3,3,3,3,3,2,3,3,1,3,3,2,1,3,2,3,1,2,3,3,1,1,1,3,1,2,1,1,2,3,2,3,1,3,1,3,3,3,1,3,3,2,2,1,3,3,3,2,1,1,1,3,1,3)

$startRow = 2
$endRow = $startRow + $values.Count - 1

# Build a 2D (N x 1) array for a single bulk range write
$data = New-Object 'object[,]' $values.Count,1
for ($i = 0; $i -lt $values.Count; $i++) {
    $data[$i,0] = $values[$i]
}

$rangeAddress = "A" + $startRow + ":A" + $endRow
$range = $ws.Range($rangeAddress)
$range.Value = $data
